# "maj template comment à la fin"
#
# The "Comment" column (currently column J, header + the 3 descriptor rows
# below it) is moved to become the last column of the header block (J:M),
# i.e. it ends up in column M. The WaitingTime / Temperature / Result
# columns (currently K, L, M) each shift one column to the left (becoming
# J, K, L respectively). Only rows 1-4 (header + the 3 metadata rows) in
# columns J:M are affected; row 5 is blank there both before and after.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the original "Comment" column's values (rows 1-4) before they're
# shifted away.
$commentVals = @(
    $ws.Range("J1").Value2,
    $ws.Range("J2").Value2,
    $ws.Range("J3").Value2,
    $ws.Range("J4").Value2
)

# Delete the J1:J4 cells, shifting K1:M4 one column to the left (into
# J1:L4). This keeps each moved cell's original type/format, it's just a
# plain column-left shift of a range (not a whole-column delete), so it
# doesn't disturb column widths/metadata.
$ws.Range("J1:J4").Delete(-4159) | Out-Null   # -4159 = xlShiftToLeft / xlToLeft

# Write the original Comment values into the now-vacated M1:M4 (the new
# last column of the header block).
for ($i = 0; $i -lt 4; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 13).Value2 = $commentVals[$i]
}
